$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new header columns (AD:AF) reusing the same header formatting
# (bold, centered, top-aligned, thin border) already applied to AC1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (Wins/Losses/Ties) for every player row.
for ($r = 2; $r -le 68; $r++) {
    $ws.Cells.Item($r, 30).Value = 73
    $ws.Cells.Item($r, 31).Value = 89
    $ws.Cells.Item($r, 32).Value = 0
}
